$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 489; existing rows 489..607 shift down to 490..608
$ws.Rows.Item(489).Insert()

# Populate the newly inserted row with the new data record
$ws.Cells.Item(489, 1).Value = 5
$ws.Cells.Item(489, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(489, 3).Value = "Maule"
$ws.Cells.Item(489, 4).Value = 45204
$ws.Cells.Item(489, 5).Value = 7
$ws.Cells.Item(489, 6).Value = 100114013
$ws.Cells.Item(489, 7).Value = "Zanahoria"
$ws.Cells.Item(489, 8).Value = "Sin especificar"
$ws.Cells.Item(489, 9).Value = "Primera"
$ws.Cells.Item(489, 10).Value = 600
$ws.Cells.Item(489, 11).Value = 5000
$ws.Cells.Item(489, 12).Value = 5500
$ws.Cells.Item(489, 13).Value = 5250
$ws.Cells.Item(489, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(489, 15).Value = "Región de Ñuble"
$ws.Cells.Item(489, 16).Value = 262
$ws.Cells.Item(489, 17).Value = 20
$ws.Cells.Item(489, 18).Value = "Hortaliza"
